# Updates values in need_to_buy.xlsx sheet: shift the daily forecast window forward by one day
# (drop oldest day, append a new day's data), per updated R-generated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 46001, 11343.5322924542, 10460.4695124852, 17376.26, 6818.90251379607, -4.03699890494666)
    ,@(3, 46002, 11416.7795921449, 10476.5256690075, 11232.26, 7118.35537109185, 265.10921000414)
    ,@(4, 46003, 11343.4491023742, 9818.75505912318, 11232.26, 7073.33537722173, 235.826268181038)
    ,@(5, 46004, 4268.510027217, 6816.60931794574, 11232.26, 6665.89411825233, 93.7601431749195)
    ,@(6, 46005, 3830.1203597031, 6839.90489998702, 11232.26, 6326.50526746956, 80.5895903106911)
    ,@(7, 46006, 9746.03784374469, 10053.6512143296, 11232.26, 8042.25692973142, 285.985339335875)
    ,@(8, 46007, 9746.03784374469, 10004.735813162, 11232.26, 8042.25692973142, 283.94719762056)
    ,@(9, 46008, 9746.03784374469, 9833.25104554299, 11232.26, 8042.25692973142, 276.801998969767)
    ,@(10, 46009, 9746.03784374469, 9877.52336463519, 11232.26, 8042.25692973142, 278.646678931942)
    ,@(11, 46010, 9746.03784374469, 9163.11506649034, 11232.26, 8042.25692973142, 248.87966650924)
    ,@(12, 46011, 8560.57524882407, 8890.22907560853, 11232.26, 7650.31988457303, 221.178706674232)
    ,@(13, 46012, 8469.04674334209, 8652.47939114308, 11232.26, 7641.92186846689, 210.922552483749)
    ,@(14, 46013, 9566.18303866457, 9394.78806183721, 11232.26, 8113.13153577614, 261.485816567223)
    ,@(15, 46014, 9566.18303866457, 9658.81213169996, 11232.26, 8113.13153577614, 266.39239279802)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = [double]$row[1]
    $ws.Cells.Item($r, 2).Value = [double]$row[2]
    $ws.Cells.Item($r, 3).Value = [double]$row[3]
    $ws.Cells.Item($r, 4).Value = [double]$row[4]
    $ws.Cells.Item($r, 5).Value = [double]$row[5]
    $ws.Cells.Item($r, 6).Value = [double]$row[6]
}

$wb.Save()
